$d = $word.ActiveDocument

# Pull the raw package XML (document.xml content is inside this) so we can do
# precise, surgical text-level edits that mirror the authoring tool's output
# (renumbered w:bookmarkStart/w:bookmarkEnd ids, a relocated _GoBack bookmark,
# and a run merge) rather than relying on the high-level object model, which
# does not expose raw bookmark ids or body-level (out-of-paragraph) bookmarks.
$xml = $d.WordOpenXML

# 1) Every existing bookmark (w:bookmarkStart/w:bookmarkEnd) gets its w:id
#    shifted up by one, because a brand-new bookmark (id 0) is inserted at
#    the very start of the body. Walk the ids from highest to lowest so a
#    freshly written id is never re-matched by a later, lower-numbered
#    replacement.
$xml = $xml -replace 'w:id="45"', 'w:id="46"'
$xml = $xml -replace 'w:id="44"', 'w:id="45"'
$xml = $xml -replace 'w:id="43"', 'w:id="44"'
$xml = $xml -replace 'w:id="42"', 'w:id="43"'
$xml = $xml -replace 'w:id="41"', 'w:id="42"'
$xml = $xml -replace 'w:id="40"', 'w:id="41"'
$xml = $xml -replace 'w:id="39"', 'w:id="40"'
$xml = $xml -replace 'w:id="38"', 'w:id="39"'
$xml = $xml -replace 'w:id="37"', 'w:id="38"'
$xml = $xml -replace 'w:id="36"', 'w:id="37"'
$xml = $xml -replace 'w:id="35"', 'w:id="36"'
$xml = $xml -replace 'w:id="34"', 'w:id="35"'
$xml = $xml -replace 'w:id="33"', 'w:id="34"'
$xml = $xml -replace 'w:id="32"', 'w:id="33"'
$xml = $xml -replace 'w:id="31"', 'w:id="32"'
$xml = $xml -replace 'w:id="30"', 'w:id="31"'
$xml = $xml -replace 'w:id="29"', 'w:id="30"'
$xml = $xml -replace 'w:id="28"', 'w:id="29"'
$xml = $xml -replace 'w:id="27"', 'w:id="28"'
$xml = $xml -replace 'w:id="26"', 'w:id="27"'
$xml = $xml -replace 'w:id="25"', 'w:id="26"'
$xml = $xml -replace 'w:id="24"', 'w:id="25"'
$xml = $xml -replace 'w:id="23"', 'w:id="24"'
$xml = $xml -replace 'w:id="22"', 'w:id="23"'
$xml = $xml -replace 'w:id="21"', 'w:id="22"'
$xml = $xml -replace 'w:id="20"', 'w:id="21"'
$xml = $xml -replace 'w:id="19"', 'w:id="20"'
$xml = $xml -replace 'w:id="18"', 'w:id="19"'
$xml = $xml -replace 'w:id="17"', 'w:id="18"'
$xml = $xml -replace 'w:id="16"', 'w:id="17"'
$xml = $xml -replace 'w:id="15"', 'w:id="16"'
$xml = $xml -replace 'w:id="14"', 'w:id="15"'
$xml = $xml -replace 'w:id="13"', 'w:id="14"'
$xml = $xml -replace 'w:id="12"', 'w:id="13"'
$xml = $xml -replace 'w:id="11"', 'w:id="12"'
$xml = $xml -replace 'w:id="10"', 'w:id="11"'
$xml = $xml -replace 'w:id="9"', 'w:id="10"'
$xml = $xml -replace 'w:id="8"', 'w:id="9"'
$xml = $xml -replace 'w:id="7"', 'w:id="8"'
$xml = $xml -replace 'w:id="6"', 'w:id="7"'
$xml = $xml -replace 'w:id="5"', 'w:id="6"'
$xml = $xml -replace 'w:id="4"', 'w:id="5"'
$xml = $xml -replace 'w:id="3"', 'w:id="4"'
$xml = $xml -replace 'w:id="2"', 'w:id="3"'
$xml = $xml -replace 'w:id="1"', 'w:id="2"'
$xml = $xml -replace 'w:id="0"', 'w:id="1"'

# 2) Insert the new _GoBack bookmark (collapsed, id 0) as the very first
#    content of the body, before the cover-page content control, exactly
#    like Word emits when the last tracked edit position is displaced by a
#    surrounding custom-XML/content-control boundary.
$xml = $xml -replace '<w:body>', '<w:body><w:bookmarkStart w:id="0" w:name="_GoBack" w:displacedByCustomXml="next"/><w:bookmarkEnd w:id="0" w:displacedByCustomXml="next"/>'

# 3) Remove the old inline _GoBack bookmark pair (now id 46) from the
#    revision table's "11:40 CEST" cell and merge the two runs it used to
#    split into a single run.
$xml = $xml -replace '<w:r><w:t xml:space="preserve">11:40 </w:t></w:r><w:bookmarkStart w:id="46" w:name="_GoBack"/><w:r><w:t>CEST</w:t></w:r><w:bookmarkEnd w:id="46"/>', '<w:r><w:t>11:40 CEST</w:t></w:r>'

$d.WordOpenXML = $xml
